# Add Game/GameMod/Championship controller tasks to the "Tasks" sheet,
# mark task #10 (row 11) as Closed, and leave the selection on F13 -
# matching the upstream commit "Add some functionality for Game controller."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- Row 11 (task #10): flip Status from Opened -> Closed ------------------
# Row 5 already carries the "Closed" style (green fill), copy its E-column
# formatting onto E11 so the cell picks up the same style index, then set
# the text.
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = "Closed"

# --- Rows 12-14: fill in the previously-empty task rows ---------------------
# Row 6 has the formatting (Task / Task Location / Priority / Status style
# indices) that all three new rows need, so copy its A:F formats down first.
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A12:F12").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:F13").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:F14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 12 - Championship controller task
$ws.Range("B12").Value = "Реализовать добавление чемпионатов в определенную лигу. Редатирование, удаление."
$ws.Range("C12").Value = "ChampionshipController`nAdd`nEdit"
$ws.Range("D12").Value = "Hight"
$ws.Range("E12").Value = "Opened"
$ws.Rows.Item(12).RowHeight = 45

# Row 13 - Game controller task
$ws.Range("B13").Value = "Реализовать добавление игр на сайт"
$ws.Range("C13").Value = "GameController`nAdd`nEdit`nDelete"
$ws.Range("D13").Value = "Hight"
$ws.Range("E13").Value = "Opened"
$ws.Rows.Item(13).RowHeight = 60

# Row 14 - GameMod controller task
$ws.Range("B14").Value = "Реализовать добавление модов на сайт"
$ws.Range("C14").Value = "GameModController`nAdd`nEdit`nDelete"
$ws.Range("D14").Value = "Hight"
$ws.Range("E14").Value = "Opened"
$ws.Rows.Item(14).RowHeight = 60

# --- Leave selection where the author left it (cell F13) -------------------
$ws.Activate()
$ws.Range("F13").Select() | Out-Null
